# 10Th - MB for single stock and added new group
#
# The sheet tracks weekly "MarketBeat rank" snapshots for a single stock.
# Each week a new pair of columns (rating-detail + plain "UN" marker) is
# inserted at the left of the history, pushing the older weekly columns to
# the right. This edit adds the snapshots for Jun_26 and Jun_27, records a
# new upgrade ("6/22/2018,Upgrades,Overweight,$25.00") for Piper Jaffray
# Companies (row 10) highlighted in the two newest columns, and appends a
# new analyst group (Benchmark, Evercore ISI) at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert three new columns (B,C,D) in front of the existing weekly
#        history (old B..E -> new E..H). ---
$ws.Range("B1:D27").EntireColumn.Insert()

# --- 2. New header row values for the inserted weekly columns. ---
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# --- 3. Populate the new B,C,D columns for every existing data row with the
#        default "UN" marker (matches the rest of the table). ---
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# --- 4. Row 10 (Piper Jaffray Companies) got a fresh upgrade on 6/22/2018,
#        recorded (and highlighted) in both new rating columns C & D. ---
$ws.Range("C10").Value = "6/22/2018,Upgrades,Overweight,`$25.00"
$ws.Range("D10").Value = "6/22/2018,Upgrades,Overweight,`$25.00"
$ws.Range("C10").Interior.ColorIndex = 42
$ws.Range("D10").Interior.ColorIndex = 42

# --- 5. New analyst group appended at the bottom of the table. ---
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
